# "Generate Report for Handback"
#
# The CI job that produces localization-status.xlsx re-ran after a handback
# completed: the zh-cn / de-de rows move from "Ready for handoff" to
# "Handed back: in sync with en-US", each row's "Latest Target File" /
# "Latest Handback File" / "Latest Handback DateTime" columns get populated,
# and a few columns on the per-language sheets are widened so the new
# (longer) values aren't truncated.

$wb = $excel.ActiveWorkbook

$mdFile  = "4eb5e790-29b6-4400-9b18-a3684345267b.md"
$mdUrl   = "https://github.com/OpenLocalizationTestOrg/oltest/blob/d8f8308abee1d913115a6466bbfebdded6be139f/e2e/4eb5e790-29b6-4400-9b18-a3684345267b.md"
$status  = "Handed back: in sync with en-US"

# Cornflowerblue (FF6495ED), matching the workbook's existing "HyperLink"
# cell style (stored by Excel as a BGR-ordered long: B*65536 + G*256 + R).
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------
# Overview sheet: widen the per-language status columns (E, F) so the
# longer status string fits.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.1675
$wsOverview.Columns.Item(6).ColumnWidth = 29.1675

# ---------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de): mark handback complete and fill in
# the target/handback file + datetime columns.
# ---------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; Xlf = "4eb5e790-29b6-4400-9b18-a3684345267b.4421b21dab6312a3e2c0564440891d670cbffa4d.zh-cn.xlf"; HandbackDate = "2016-08-12 07:12:53" },
    @{ Name = "de-de"; Xlf = "4eb5e790-29b6-4400-9b18-a3684345267b.4421b21dab6312a3e2c0564440891d670cbffa4d.de-de.xlf"; HandbackDate = "2016-08-12 07:13:06" }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Status column (C) -> handed back
    $ws.Range("C2").Value = $status

    # Widen Status (C), Latest Target File (I) and Latest Handback File (J)
    $ws.Columns.Item(3).ColumnWidth = 29.1675
    $ws.Columns.Item(9).ColumnWidth = 39.1667
    $ws.Columns.Item(10).ColumnWidth = 39.1667

    # Latest Target File (I2): the handed-back markdown file, hyperlinked
    # just like the source file name in column A.
    $ws.Range("I2").Value = $mdFile
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdUrl, "", "", $mdFile)
    $ws.Range("I2").Font.Underline = $true
    $ws.Range("I2").Font.Color = $hyperlinkColor

    # Latest Handback File (J2): the xliff that was handed back.
    $ws.Range("J2").Value = $lang.Xlf

    # Latest Handback DateTime (K2).
    $ws.Range("K2").Value = $lang.HandbackDate
}
